# The sheet contains a weekly price list for "Espinaca" (Vega Modelo de
# Temuco). A new weekly record was inserted as row 70, pushing the
# previously-existing rows 70..101 down to 71..102 (all of their data is
# unchanged by this edit - only their row position moves).
#
# Insert a new row at position 70 (Excel shifts rows 70-101 down to 71-102,
# carrying over formatting/styles just like typing a new row in the UI),
# then populate the new row 70 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above current row 70; this shifts existing rows
# 70..101 down to 71..102.
$ws.Rows.Item(70).Insert()

# Fill in the new row 70 with the new record. Columns A, B, C, E, F, G, H,
# I, N, O, Q, R follow the same fixed template used by every other row in
# this sheet (same market/product/region/unit/classification); only D
# (Fecha), J (Volumen), K/L/M (prices) and P (Precio $/Kg) vary per row.
$ws.Cells.Item(70, 1).Value() = 10
$ws.Cells.Item(70, 2).Value() = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value() = "La Araucanía"
$ws.Cells.Item(70, 4).Value() = 44518
$ws.Cells.Item(70, 5).Value() = 9
$ws.Cells.Item(70, 6).Value() = 100112012
$ws.Cells.Item(70, 7).Value() = "Espinaca"
$ws.Cells.Item(70, 8).Value() = "Sin especificar"
$ws.Cells.Item(70, 9).Value() = "Primera"
$ws.Cells.Item(70, 10).Value() = 50
$ws.Cells.Item(70, 11).Value() = 8000
$ws.Cells.Item(70, 12).Value() = 8000
$ws.Cells.Item(70, 13).Value() = 8000
$ws.Cells.Item(70, 14).Value() = "$/docena de atados"
$ws.Cells.Item(70, 15).Value() = "Región de La Araucanía"
$ws.Cells.Item(70, 16).Value() = 2667
$ws.Cells.Item(70, 17).Value() = 3
$ws.Cells.Item(70, 18).Value() = "Hortaliza"
